$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testing")
$ws.Activate()

# Fix existing rows 2-8: clamp_angle (column J) changes from 30 to 0
for ($r = 2; $r -le 8; $r++) {
    $ws.Range("J$r").Value = 0
}

# Append new trial rows 9-25, following the same pattern as the existing rows
# (trial_num increments, target_amp = 3, terminal_feedback = 1, everything else 0)
for ($r = 9; $r -le 25; $r++) {
    $trial = $r - 1
    $ws.Range("A$r").Value = $trial
    $ws.Range("B$r").Value = 0
    $ws.Range("C$r").Value = 3
    $ws.Range("D$r").Value = 0
    $ws.Range("E$r").Value = 1
    $ws.Range("F$r").Value = 0
    $ws.Range("G$r").Value = 0
    $ws.Range("H$r").Value = 0
    $ws.Range("I$r").Value = 0
    $ws.Range("J$r").Value = 0
    $ws.Range("K$r").Value = 0
    $ws.Range("L$r").Value = 0
}

# Match formatting of the new rows to the existing data rows above them
$ws.Range("A8:L8").Copy()
$ws.Range("A9:L25").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Restore the last active-cell selection on this sheet
$ws.Range("Q12").Select()
